$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/sheet tab
$ws.Name = "ShearF"

# Add new row 16 with the next Gaussian Quadrature scheme results
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 0.985526333062672
$ws.Range("D16").Value = 1.045570974701326
$ws.Range("E16").Value = 0.9877485424670935
$ws.Range("F16").Value = 0.985526333062672
$ws.Range("G16").Value = 1.026035405532553
$ws.Range("H16").Value = 0.9671672344821337
$ws.Range("I16").Value = 0.9877192943386085
$ws.Range("J16").Value = 1.045570974701326
$ws.Range("K16").Value = 1.01665975858421
$ws.Range("L16").Value = 1.001093045823441
$ws.Range("M16").Value = 0.9999612974307311

# Copy styling from row 15 (A15 has bold/centered/border style) to A16
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
